$d = $word.ActiveDocument

# --- 1. Title: "Visionsdokument for Flextur" -> append "-system" as a new run ---
$p1 = $d.Paragraphs.Item(1)
$endPos = $p1.Range.End
$rTitle = $d.Range($endPos - 1, $endPos - 1)
$rTitle.InsertAfter("-system")

# --- 2. "bestillingsprocessen." paragraph: remove "sindsro " and change the ending ---
$d.Content.Find.Execute("med sindsro ikke være i tvivl om, at kørslen er bestilt. Endvidere", $true, $false, $false, $false, $false, $true, 1, $false, "med ikke være i tvivl om, at hvordan man vha. systemet bestiller kørsel. Endvidere", 2) | Out-Null

# --- 3. "De specifikke krav..." paragraph ---
$d.Content.Find.Execute("af systemet. Her skal brugeren", $true, $false, $false, $false, $false, $true, 1, $false, "af systemet. Hvis kunden benytter systemet for første gang skal vedkommende kunne oprette sig som bruger af systemet, Herefter skal brugeren", 2) | Out-Null

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text.StartsWith("De specifikke krav")) {
        $specIdx = $i
        break
    }
}
$pSpec = $d.Paragraphs.Item($specIdx)
$specEnd = $pSpec.Range.End
$rSpec = $d.Range($specEnd - 1, $specEnd - 1)
$rSpec.InsertAfter("Endvidere er en brugeroversigt over tidligere kørsler ønskværdigt.")

# --- 4. Insert new paragraph before "Medarbejderens interesser..." paragraph ---
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text.StartsWith("Medarbejderens interesser ligger")) {
        $medIdx = $i
        break
    }
}
$pMed = $d.Paragraphs.Item($medIdx)
$pMed.Range.InsertParagraphBefore()
$pNewMed = $d.Paragraphs.Item($medIdx)
$pNewMed.Range.Text = "Som den anden primære bruger af systemet er det vigtigt for medarbejderen, at vedkommende let kan kende forskel på kunde- og medarbejderinterfacet i systemet. Medarbejderens tilgangsvinkel til systemet skal være målrettet til distributionsbrug og så brugervenligt og tilpasset som muligt. "

# --- 5. Rewrite "Medarbejderens interesser ligger i," -> "Medarbejderens specifikke krav er," ---
$d.Content.Find.Execute("Medarbejderens interesser ligger i, at kundens kørselsbestilling", $true, $false, $false, $false, $false, $true, 1, $false, "Medarbejderens specifikke krav er, at kundens kørselsbestilling", 2) | Out-Null

# --- 6. Append new sentence after "... kan tildeles en bil. Endvidere " ---
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text.StartsWith("Medarbejderens specifikke krav er") -or $d.Paragraphs.Item($i).Range.Text.StartsWith("Medarbejderens interesser ligger")) {
        $medReqIdx = $i
        break
    }
}
$pMedReq = $d.Paragraphs.Item($medReqIdx)
$medReqEnd = $pMedReq.Range.End
$rMedReq = $d.Range($medReqEnd - 1, $medReqEnd - 1)
$rMedReq.InsertAfter("ønsker medarbejderen at få vist en oversigt over afholdte ture hhv. pr. bruger og pr. tidsinterval.")

# --- 7. "Midttrafiks interesser..." paragraph rewrite ---
$d.Content.Find.Execute("fejl eller uoverensstemmelser mellem bruger og Flextur-produktet, der er på baggrund af systemet. ", $true, $false, $false, $false, $false, $true, 1, $false, "fejl , snyd eller uoverensstemmelser mellem bruger og Flextur-produktet, der skyldes systemet.", 2) | Out-Null

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text.StartsWith("Midttrafiks interesser")) {
        $midInterestIdx = $i
        break
    }
}
$pMidInterest = $d.Paragraphs.Item($midInterestIdx)
$midInterestEnd = $pMidInterest.Range.End
$rMidInterest = $d.Range($midInterestEnd - 1, $midInterestEnd - 1)
$rMidInterest.InsertAfter(" Endvidere skal kørselsbestillingen og kørselsgodkendelsesprocessen ske så effektivt og minimalistisk som muligt, så kun de nødvendige informationer bliver delt kunden og medarbejderne imellem. Endvidere er det i firmaets store interesse, at systemet optimeres til at give alle brugere af systemet den bedste brugeroplevelse, hvor kunden får den mest optimale kørsel til den korrekte pris.")

# --- 8. Delete the trailing empty Heading2 + 3 empty paragraphs after the Midttrafik paragraph ---
# Find current paragraph index of "Midttrafiks interesser..." paragraph, then remove the next 4 paragraphs.
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text.StartsWith("Midttrafiks interesser")) {
        $midIdx = $i
        break
    }
}
$delStart = $d.Paragraphs.Item($midIdx + 1).Range.Start
$delEnd = $d.Paragraphs.Item($midIdx + 4).Range.End
$rDel = $d.Range($delStart, $delEnd)
$rDel.Delete()

# --- 9. Feature list item rewrites ---
$d.Content.Find.Execute("Systemet udfører sletning af kunde", $true, $false, $false, $false, $false, $true, 1, $false, "Sletning af kunde", 2) | Out-Null
$d.Content.Find.Execute("Systemet udfører redigering af kunde", $true, $false, $false, $false, $false, $true, 1, $false, "Redigering af kunde", 2) | Out-Null
$d.Content.Find.Execute("Kørselsoverblik", $true, $false, $false, $false, $false, $true, 1, $false, "Kørselsadministration for bruger", 2) | Out-Null
$d.Content.Find.Execute("Systemet udfører bestilling af kørsel", $true, $false, $false, $false, $false, $true, 1, $false, "Bestilling af kørsel", 2) | Out-Null
$d.Content.Find.Execute("Systemet udfører sletning af kørsel", $true, $false, $false, $false, $false, $true, 1, $false, "Visning af kørselshistorik", 2) | Out-Null
$d.Content.Find.Execute("Registerindsigt", $true, $false, $false, $false, $false, $true, 1, $false, "Kørselsadministration for Midttrafik", 2) | Out-Null
$d.Content.Find.Execute("Systemet udfører visning af brugerhistorik", $true, $false, $false, $false, $false, $true, 1, $false, "Godkendelse af kørsel", 2) | Out-Null
$d.Content.Find.Execute("Systemet udfører visning af fremtidige kørsler", $true, $false, $false, $false, $false, $true, 1, $false, "Tildeling af bil til kørsel", 2) | Out-Null
$d.Content.Find.Execute("Systemet udfører visning af tidligere kørsler", $true, $false, $false, $false, $false, $true, 1, $false, "Visning af fremtidige kørsler med pris", 2) | Out-Null

# --- 10. Add new list paragraph "Visning af tidligere kørsler pr. bruger og tidsinterval" after "Visning af fremtidige kørsler med pris" ---
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text.StartsWith("Visning af fremtidige kørsler med pris")) {
        $lastListIdx = $i
        break
    }
}
$pLastList = $d.Paragraphs.Item($lastListIdx)
$pLastList.Range.InsertParagraphAfter()
$pNewList = $d.Paragraphs.Item($lastListIdx + 1)
$pNewList.Range.Text = "Visning af tidligere kørsler pr. bruger og tidsinterval"
